# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-67) from 2023-09-17 (serial 45186) to 2023-09-19 (serial 45188).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 67; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
